# Auto-update draw results: append the 2025-09-28 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 12

# All columns in this sheet hold text (dates, phase codes and the "HH:MM..."
# style draw code are stored as strings, not numbers/dates), so force the
# new row's number format to Text before writing the values - otherwise the
# host would auto-coerce "2025-09-28" to a date serial and "250928" to a
# plain number.
$ws.Range("A" + $newRow + ":E" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-09-28"
$ws.Range("B" + $newRow).Value = "Pick 4"
$ws.Range("C" + $newRow).Value = "250928"
$ws.Range("D" + $newRow).Value = "3-4-6-1"
$ws.Range("E" + $newRow).Value = "2025-09-28T21:34:50.531+04:00"
